$wb = $excel.ActiveWorkbook

# --- Shared values used across sheets -------------------------------------------------
$handbackFileDisplay = "b28016bb-3dfa-4e58-bb39-f5504b165943.md"
$handbackFileUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/00a0dc43d75b1314d6c45092b2bedb2eb9198a78/e2e/b28016bb-3dfa-4e58-bb39-f5504b165943.md"
$versionMismatchMsg  = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f926f97d9702c07fba265930bce352745a7673ad/e2e/b28016bb-3dfa-4e58-bb39-f5504b165943.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/00a0dc43d75b1314d6c45092b2bedb2eb9198a78/e2e/b28016bb-3dfa-4e58-bb39-f5504b165943.md."

# Hyperlink cell font (matches the existing "HyperLink" style already used in the workbook)
$hyperlinkColor = 15570276   # RGB(0x64, 0x95, 0xED) -> matches font color FF6495ED

# --- zh-cn sheet (row 7 : b28016bb-3dfa-4e58-bb39-f5504b165943) ------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("K7").Value = "2016-09-05 07:03:47"
$wsZhCn.Range("P7").Value = $versionMismatchMsg

$wsZhCn.Range("J7").Value = "b28016bb-3dfa-4e58-bb39-f5504b165943.34e4728ec4160f0fc86791637e56f87f88c357a1.zh-cn.xlf"

$zhCnLink = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), $handbackFileUrl, "", "", $handbackFileDisplay)
$wsZhCn.Range("I7").Font.Name = "Calibri"
$wsZhCn.Range("I7").Font.Underline = $true
$wsZhCn.Range("I7").Font.Color = $hyperlinkColor

# --- de-de sheet (row 7 : b28016bb-3dfa-4e58-bb39-f5504b165943) ------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("K7").Value = "2016-09-05 07:03:55"
$wsDeDe.Range("P7").Value = $versionMismatchMsg

$wsDeDe.Range("J7").Value = "b28016bb-3dfa-4e58-bb39-f5504b165943.34e4728ec4160f0fc86791637e56f87f88c357a1.de-de.xlf"

$deDeLink = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), $handbackFileUrl, "", "", $handbackFileDisplay)
$wsDeDe.Range("I7").Font.Name = "Calibri"
$wsDeDe.Range("I7").Font.Underline = $true
$wsDeDe.Range("I7").Font.Color = $hyperlinkColor
